$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosGenerales")

# Set B6 to the text value "25" (cell is text-formatted, so keep it as a string)
$ws.Range("B6").Value = "25"

# Move the active selection from B5 to B6
$ws.Activate()
$ws.Range("B6").Select()
